$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 3-25: Runmode flips from Y to N
$ws.Range("D3:D25").Value = "N"

# Row 2 result flips from SKIP to PASS (test case now passes)
$ws.Range("E2").Value = "PASS"

# Rows 16-25: TestCase_B25/B26/B27 duplicates removed, remaining
# TestCase_B18..B27 renumbered down to TestCase_B15..B24
$ws.Range("A16").Value = "TestCase_B15"
$ws.Range("A17").Value = "TestCase_B16"
$ws.Range("A18").Value = "TestCase_B17"
$ws.Range("A19").Value = "TestCase_B18"
$ws.Range("A20").Value = "TestCase_B19"
$ws.Range("A21").Value = "TestCase_B20"
$ws.Range("A22").Value = "TestCase_B21"
$ws.Range("A23").Value = "TestCase_B22"
$ws.Range("A24").Value = "TestCase_B23"
$ws.Range("A25").Value = "TestCase_B24"

# Move the active selection to C16, matching the saved cursor position
$null = $ws.Range("C16").Select()
